# #5: property aircraft done
# Fix the "property_category" column values that were mistakenly left as
# "land" on the 建物 (Building) and 汽車 (Car) sheets.

$wb = $excel.ActiveWorkbook

# 建物 (Building) sheet: column I is "property_category", rows 2-9 should be "building"
$wsBuilding = $wb.Worksheets.Item("建物")
for ($r = 2; $r -le 9; $r++) {
    $wsBuilding.Cells.Item($r, 9).Value = "building"
}

# 汽車 (Car) sheet: column H is "property_category", rows 2-3 should be "car"
$wsCar = $wb.Worksheets.Item("汽車")
for ($r = 2; $r -le 3; $r++) {
    $wsCar.Cells.Item($r, 8).Value = "car"
}
